$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.087.13'
$ws.Range("E2").Value = '  -5.63%  '

$ws.Range("D3").Value = '3.335.64'
$ws.Range("E3").Value = '  -2.58%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '564.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.48%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.33'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.01%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").Value = '3.335.49'
$ws.Range("E8").Value = '  -2.61%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.472'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.49%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.42'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.25%  '

$ws.Range("E11").Value = '  -4.67%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.376'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.34%  '

$ws.Range("D13").Value = '3.905.34'
$ws.Range("E13").Value = '  -2.50%  '

$ws.Range("E14").Value = '  -0.23%  '

$ws.Range("D15").Value = '3.344.06'
$ws.Range("E15").Value = '  -2.24%  '

$ws.Range("E16").Value = '  -4.09%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '24.62'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.16%  '

$ws.Range("D18").Value = '60.179.16'
$ws.Range("E18").Value = '  -5.51%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.69'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.93%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.47'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.34%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.11'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -7.37%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '354.41'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -7.29%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.558'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.81%  '

$ws.Range("D24").Value = '3.469.45'
$ws.Range("E24").Value = '  -2.60%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.22%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '69.28'
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000111'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.95%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.66'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +18.71%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.47'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.45%  '

$ws.Range("E30").Value = '  -0.11%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.95'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.60%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.153'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.62%  '

$ws.Range("E33").Value = '  -3.03%  '

$ws.Range("E34").Value = '  -0.05%  '

$ws.Range("D35").Value = '3.365.48'
$ws.Range("E35").Value = '  -2.60%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '22.92'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.04%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.38'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.79%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.89'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.58%  '

$ws.Range("E39").Value = '  +0.47%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '159.18'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.00%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0768'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.15%  '

$ws.Range("E42").Value = '  +0.05%  '

$ws.Range("E43").Value = '  +2.12%  '

$ws.Range("E44").Value = '  -4.30%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '40.78'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.48%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.18'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +7.96%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.57'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.23%  '

$ws.Range("E48").Value = '  -1.34%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.80'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.62%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.35'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +10.41%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.890'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.47%  '
